$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D20").Value = "Damage SFX (x2)"
$ws.Range("E20").Value = "Completed"

$ws.Range("E26").Value = "Completed"
$ws.Range("E27").Value = "Completed"

$ws.Range("A28").Value = "UIDrag"
$ws.Range("B28").Value = "Sound of a slider being dragged"
$ws.Range("C28").Value = "Interface"
$ws.Range("D28").Value = "Drag SFX (x2)"
$ws.Range("E28").Value = "Completed"

# Match the author's final selection position in the saved view
$ws.Range("E21").Select() | Out-Null
